$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row cell updates scraped from the commit diff. Columns B/C hold plain
# text (coin name / link) and are safe to assign directly. Columns D/E/G
# hold numeric-looking text (price, % change, "hour" counter) that Excel's
# COM layer would otherwise auto-coerce into a Number -- force them back to
# Text via NumberFormat "@" before the write, then ClearFormats so the
# cell's style id is left untouched (matches the source workbook, which
# carries no explicit style on these cells).
$updates = @(
    @{ Row=2; D='257.42'; E='-0.46%'; G='5' }
    @{ Row=3; D='27.03'; E='0.56%'; G='5' }
    @{ Row=4; D='4.565'; E='-6.07%'; G='5' }
    @{ Row=5; D='0.05899'; G='5' }
    @{ Row=6; D='6.630'; E='-0.86%'; G='5' }
    @{ Row=7; D='0.8529'; E='-2.74%'; G='5' }
    @{ Row=8; D='0.9413'; E='-1.76%'; G='5' }
    @{ Row=9; B='One'; C='https://coinranking.com/coin/6Lga5NiXX3rT+one-one'; D='0.01038'; E='1,607.74%'; G='5' }
    @{ Row=10; B='WazirX'; C='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; D='0.1391'; E='-1.69%'; G='5' }
    @{ Row=11; B='LiechtensteinCryptoassetsExchange'; C='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; D='0.05073'; E='40.86%'; G='5' }
    @{ Row=12; B='MandalaExchangeToken'; C='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; D='0.07083'; E='-1.94%'; G='5' }
    @{ Row=13; B='BitrueCoin'; C='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; D='0.03069'; E='-2.48%'; G='5' }
    @{ Row=14; B='BitMartToken'; C='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; D='0.09118'; E='-1.28%'; G='5' }
    @{ Row=15; B='BitForexToken'; C='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; D='0.001526'; E='-0.79%'; G='5' }
    @{ Row=16; D='0.006121'; E='2.87%'; G='5' }
    @{ Row=17; D='3.491'; E='0.19%'; G='5' }
    @{ Row=18; D='3.179'; E='-1.31%'; G='5' }
    @{ Row=19; G='5' }
    @{ Row=20; D='0.3052'; E='-2.91%'; G='5' }
    @{ Row=21; D='0.1269'; E='-2.86%'; G='5' }
    @{ Row=22; D='3.942'; E='11.94%'; G='5' }
    @{ Row=23; D='0.04253'; E='1.02%'; G='5' }
    @{ Row=24; D='0.001219'; E='-0.30%'; G='5' }
    @{ Row=25; D='0.004286'; E='-5.12%'; G='5' }
    @{ Row=26; E='0.05%'; G='5' }
    @{ Row=27; E='2.08%'; G='5' }
    @{ Row=28; G='5' }
    @{ Row=29; G='5' }
    @{ Row=30; G='5' }
    @{ Row=31; G='5' }
    @{ Row=32; G='5' }
    @{ Row=33; G='5' }
    @{ Row=34; G='5' }
    @{ Row=35; G='5' }
    @{ Row=36; G='5' }
    @{ Row=37; G='5' }
    @{ Row=38; G='5' }
    @{ Row=39; G='5' }
    @{ Row=40; E='-0.30%'; G='5' }
    @{ Row=41; B='KickToken'; C='https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'; D='0.006230'; E='5.82%'; G='5' }
    @{ Row=42; B='BKEXToken'; C='https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'; D='0.1100'; E='-0.34%'; G='5' }
    @{ Row=43; D='0.002200'; E='15.84%'; G='5' }
    @{ Row=44; D='0.01405'; E='33.95%'; G='5' }
    @{ Row=45; D='0.00005343'; E='-2.63%'; G='5' }
    @{ Row=46; E='0.05%'; G='5' }
    @{ Row=47; G='5' }
    @{ Row=48; D='0.2517'; E='11,743.10%'; G='5' }
    @{ Row=49; E='0.05%'; G='5' }
    @{ Row=50; E='0.05%'; G='5' }
    @{ Row=51; G='5' }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($u.ContainsKey('B')) {
        $ws.Cells.Item($r, 2).Value = $u.B
    }
    if ($u.ContainsKey('C')) {
        $ws.Cells.Item($r, 3).Value = $u.C
    }
    if ($u.ContainsKey('D')) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.ClearFormats()
    }
    if ($u.ContainsKey('E')) {
        $cell = $ws.Cells.Item($r, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.ClearFormats()
    }
    if ($u.ContainsKey('G')) {
        $cell = $ws.Cells.Item($r, 7)
        $cell.NumberFormat = "@"
        $cell.Value = $u.G
        $cell.ClearFormats()
    }
}
